$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 167, shifting rows 167-171 down to 168-172
$ws.Rows.Item(167).Insert()

# Fill in the new row 167 with the new data
$ws.Cells.Item(167, 1).Value = 8
$ws.Cells.Item(167, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(167, 3).Value = 'Coquimbo'
$ws.Cells.Item(167, 4).Value = 44448
$ws.Cells.Item(167, 5).Value = 4
$ws.Cells.Item(167, 6).Value = 100112032
$ws.Cells.Item(167, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(167, 8).Value = 'Sin especificar'
$ws.Cells.Item(167, 9).Value = 'Primera'
$ws.Cells.Item(167, 10).Value = 500
$ws.Cells.Item(167, 11).Value = 14000
$ws.Cells.Item(167, 12).Value = 15000
$ws.Cells.Item(167, 13).Value = 14500
$ws.Cells.Item(167, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(167, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(167, 16).Value = 290
$ws.Cells.Item(167, 17).Value = 50
$ws.Cells.Item(167, 18).Value = 'Hortaliza'

$wb.Save()
